# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-data updates to the Sheets workbook
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*)

$wb = $excel.ActiveWorkbook

# ----- Sheet "ALC" -----
$ws = $wb.Worksheets("ALC")

# Row 33
$ws.Range("H33").Value = 464.63635
$ws.Range("I33").Value = 173
$ws.Range("K33").Value = 173
$ws.Range("M33").Value = 56

# Row 86
$ws.Range("H86").Value = 5649.8335
$ws.Range("J86").Value = 5649.8335
$ws.Range("L86").Value = 5649.8335
$ws.Range("N86").Value = -7895.8335

# Row 89
$ws.Range("H89").Value = 5649.8335
$ws.Range("J89").Value = 5649.8335
$ws.Range("L89").Value = 28249.1675
$ws.Range("N89").Value = -39481.1675

# Row 97
$ws.Range("H97").Value = 4530.4443
$ws.Range("J97").Value = 4044.25
$ws.Range("L97").Value = 12132.75
$ws.Range("N97").Value = -13124.75

# Row 99
$ws.Range("H99").Value = 182.57143
$ws.Range("I99").Value = 182.57143
$ws.Range("K99").Value = 547.71429
$ws.Range("M99").Value = 950.28571

# Row 100
$ws.Range("H100").Value = 759.5
$ws.Range("I100").Value = 759.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 759.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -218.5
$ws.Range("N100").ClearContents()

# Row 111
$ws.Range("H111").Value = 1386.3334
$ws.Range("I111").Value = 1173.7646
$ws.Range("K111").Value = 3521.2938
$ws.Range("M111").Value = -454.2937999999999

# Row 138
$ws.Range("H138").Value = 5708.4
$ws.Range("J138").Value = 5168.2583
$ws.Range("L138").Value = 15504.7749
$ws.Range("N138").Value = -25784.7749


# ----- Sheet "ARM" -----
$ws = $wb.Worksheets("ARM")

# Row 32
$ws.Range("H32").Value = 18356.643
$ws.Range("I32").Value = 8170.6387
$ws.Range("K32").Value = 8170.6387
$ws.Range("M32").Value = -7883.6387

# Row 45
$ws.Range("H45").Value = 2888.2856
$ws.Range("I45").Value = 1609.5
$ws.Range("J45").Value = 3399.8
$ws.Range("K45").Value = 1609.5
$ws.Range("L45").Value = 3399.8
$ws.Range("M45").Value = -1232.5
$ws.Range("N45").Value = -4153.8

# Row 61
$ws.Range("H61").Value = 2470.3333
$ws.Range("I61").Value = 2320.4285
$ws.Range("K61").Value = 2320.4285
$ws.Range("M61").Value = -2108.4285

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 74
$ws.Range("H74").Value = 3528
$ws.Range("I74").Value = 1581.375
$ws.Range("K74").Value = 1581.375
$ws.Range("M74").Value = -707.375

# Row 77
$ws.Range("H77").Value = 3528
$ws.Range("I77").Value = 1581.375
$ws.Range("K77").Value = 7906.875
$ws.Range("M77").Value = -3538.875

# Row 97
$ws.Range("H97").Value = 592.2273
$ws.Range("J97").Value = 632.25
$ws.Range("L97").Value = 632.25
$ws.Range("N97").Value = -1624.25

# Row 136
$ws.Range("H136").Value = 2470.3333
$ws.Range("I136").Value = 2320.4285
$ws.Range("K136").Value = 6961.2855
$ws.Range("M136").Value = -4411.2855


# ----- Sheet "BSM" -----
$ws = $wb.Worksheets("BSM")

# Row 105
$ws.Range("H105").Value = 2920.6
$ws.Range("I105").Value = 2341.8
$ws.Range("K105").Value = 2341.8
$ws.Range("M105").Value = -594.8000000000002


# ----- Sheet "CRP" -----
$ws = $wb.Worksheets("CRP")

# Row 31
$ws.Range("H31").Value = 5128.2573
$ws.Range("I31").Value = 4994.6665
$ws.Range("K31").Value = 4994.6665
$ws.Range("M31").Value = -4699.6665

# Row 34
$ws.Range("H34").Value = 5128.2573
$ws.Range("I34").Value = 4994.6665
$ws.Range("K34").Value = 4994.6665
$ws.Range("M34").Value = -4792.6665

# Row 109
$ws.Range("H109").Value = 52643.09
$ws.Range("J109").Value = 52643.09
$ws.Range("L109").Value = 52643.09
$ws.Range("N109").Value = -54723.09

# Row 121
$ws.Range("H121").Value = 80000
$ws.Range("J121").Value = 80000
$ws.Range("L121").Value = 80000
$ws.Range("N121").Value = -82620

# Row 134
$ws.Range("H134").Value = 3181.3
$ws.Range("I134").Value = 2769.3076
$ws.Range("K134").Value = 8307.9228
$ws.Range("M134").Value = -5772.9228


# ----- Sheet "GSM" -----
$ws = $wb.Worksheets("GSM")

# Row 43
$ws.Range("H43").Value = 2053.4
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 80
$ws.Range("H80").Value = 6188.6665
$ws.Range("I80").Value = 3750
$ws.Range("K80").Value = 3750
$ws.Range("M80").Value = -2752

# Row 83
$ws.Range("H83").Value = 6188.6665
$ws.Range("I83").Value = 3750
$ws.Range("K83").Value = 18750
$ws.Range("M83").Value = -13758


# ----- Sheet "LTW" -----
$ws = $wb.Worksheets("LTW")

# Row 22
$ws.Range("H22").Value = 4299.6665
$ws.Range("J22").Value = 4299.6665
$ws.Range("L22").Value = 4299.6665
$ws.Range("N22").Value = -4889.6665

# Row 27
$ws.Range("H27").Value = 4299.6665
$ws.Range("J27").Value = 4299.6665
$ws.Range("L27").Value = 4299.6665
$ws.Range("N27").Value = -4513.6665

# Row 64
$ws.Range("H64").Value = 22500
$ws.Range("J64").Value = 22500
$ws.Range("L64").Value = 22500
$ws.Range("N64").Value = -22950

# Row 67
$ws.Range("H67").Value = 22500
$ws.Range("J67").Value = 22500
$ws.Range("L67").Value = 22500
$ws.Range("N67").Value = -24060

# Row 93
$ws.Range("H93").Value = 1132.4445
$ws.Range("J93").Value = 1862.5
$ws.Range("L93").Value = 1862.5
$ws.Range("N93").Value = -4358.5


# ----- Sheet "WVR" -----
$ws = $wb.Worksheets("WVR")

# Row 58
$ws.Range("H58").Value = 7540
$ws.Range("J58").Value = 10000
$ws.Range("L58").Value = 10000
$ws.Range("N58").Value = -10616

# Row 74
$ws.Range("H74").Value = 22596.834
$ws.Range("I74").Value = 3994
$ws.Range("J74").Value = 26317.4
$ws.Range("K74").Value = 3994
$ws.Range("L74").Value = 26317.4
$ws.Range("M74").Value = -3058
$ws.Range("N74").Value = -28189.4

# Row 77
$ws.Range("H77").Value = 22596.834
$ws.Range("I77").Value = 3994
$ws.Range("J77").Value = 26317.4
$ws.Range("K77").Value = 11982
$ws.Range("L77").Value = 78952.20000000001
$ws.Range("M77").Value = -7302
$ws.Range("N77").Value = -88312.20000000001

# Row 109
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774

# Row 132
$ws.Range("H132").Value = 2742.8235
$ws.Range("I132").Value = 2275.7273
$ws.Range("K132").Value = 6827.1819
$ws.Range("M132").Value = -4297.1819

# Row 141
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

